$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.614.37"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "2.485.78"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.558"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").Value = "2.518.22"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.333"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").Value = "2.930.86"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").Value = "58.461.79"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Value = "2.507.53"
$ws.Range("E18").Value = "  +2.04%  "
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "322.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.995"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.992"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").Value = "0.0₃0756"
$ws.Range("E29").Value = "  +2.28%  "
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.31%  "
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.994"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.780"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "278.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.599"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "123.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0919"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0502"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.34%  "
